$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.272.89"
$ws.Range("E2").Value = "  +3.82%  "
$ws.Range("D3").Value = "3.210.69"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'539.42"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  +4.53%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.532"
$ws.Range("E8").Value = "  +3.28%  "
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("D11").Value = "'0.433"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("D12").Value = "3.763.44"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "'26.14"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "60.317.95"
$ws.Range("E16").Value = "  +3.72%  "
$ws.Range("D17").Value = "3.196.94"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "'6.26"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'13.17"
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("E20").Value = "  +2.38%  "
$ws.Range("D21").Value = "'383.22"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "'0.531"
$ws.Range("E23").Value = "  +3.04%  "
$ws.Range("D24").Value = "'70.25"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'8.87"
$ws.Range("E25").Value = "  +11.85%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.171"
$ws.Range("E26").Value = "  +2.14%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "0.0₃0912"
$ws.Range("E28").Value = "  +3.65%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.46"
$ws.Range("E30").Value = "  +3.37%  "
$ws.Range("E31").Value = "  +5.69%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "'6.20"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("E33").Value = "  +5.43%  "
$ws.Range("D34").Value = "'6.64"
$ws.Range("E34").Value = "  +6.09%  "
$ws.Range("D35").Value = "'156.88"
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "2.784.47"
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("D38").Value = "'25.84"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "'0.0710"
$ws.Range("E39").Value = "  +4.58%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "'4.27"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "'39.87"
$ws.Range("E42").Value = "  +3.74%  "
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("E44").Value = "  +4.00%  "
$ws.Range("D45").Value = "3.251.95"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "'6.18"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "'0.809"
$ws.Range("E49").Value = "  +8.12%  "
$ws.Range("D50").Value = "'20.77"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.04%  "
